$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-style the existing data rows (7-22) so they drop the old "font 2"
#    (Arial / theme color) in favour of the default font already used by
#    rows 2-6, matching the new look-and-feel from the diff.
# ---------------------------------------------------------------------------

# Column A (timestamp cells) rows 7-22: same effective style as A2 (centered,
# vertically centered, date/time number format) but without wrap text.
$srcA = $ws.Range("A2")
$dstA = $ws.Range("A7:A22")
$srcA.Copy()
$dstA.PasteSpecial(-4122)
$dstA.WrapText = $false

# Columns B & C rows 7-12, and column F rows 7-12: becomes centered +
# vertically centered + wrap text (same as B2).
$srcB = $ws.Range("B2")

$dstBC = $ws.Range("B7:C12")
$srcB.Copy()
$dstBC.PasteSpecial(-4122)

$dstF1 = $ws.Range("F7:F12")
$srcB.Copy()
$dstF1.PasteSpecial(-4122)

# Columns D & E rows 7-22: keep left/default horizontal alignment + wrap
# text, but become vertically centered.
$dstDE = $ws.Range("D13:E22")
$srcD = $ws.Range("D7")
$srcD.Copy()
$dstDE.PasteSpecial(-4122)
$ws.Range("D7:E22").VerticalAlignment = -4108

# Columns B & C rows 13-22, and column F rows 13-22: centered + vertically
# centered, no wrap text.
$dstBC2 = $ws.Range("B13:C22")
$srcB.Copy()
$dstBC2.PasteSpecial(-4122)
$dstBC2.WrapText = $false

$dstF2 = $ws.Range("F13:F22")
$srcB.Copy()
$dstF2.PasteSpecial(-4122)
$dstF2.WrapText = $false

# ---------------------------------------------------------------------------
# 2. Append the three new survey responses (rows 23-25).
# ---------------------------------------------------------------------------

$ws.Cells.Item(23,1).Value = 44301.4706641088
$ws.Cells.Item(23,2).Value = "Não"
$ws.Cells.Item(23,3).Value = "Não"
$ws.Cells.Item(23,4).Value = "Nenhuma"
$ws.Cells.Item(23,5).Value = "Nenhuma"
$ws.Cells.Item(23,6).Value = "Não"

$ws.Cells.Item(24,1).Value = 44301.83158115741
$ws.Cells.Item(24,2).Value = "Não"
$ws.Cells.Item(24,3).Value = "Não"
$ws.Cells.Item(24,4).Value = "Monitoria em alguma disciplina Específica (obrigatória ou geral)"
$ws.Cells.Item(24,5).Value = "Projetos pessoais (websites, aplicativos próprios e etc), Nenhuma"
$ws.Cells.Item(24,6).Value = "Não"

$ws.Cells.Item(25,1).Value = 44303.43661967592
$ws.Cells.Item(25,2).Value = "Sim"
$ws.Cells.Item(25,3).Value = "Não"
$ws.Cells.Item(25,4).Value = "Participação em Eventos como Ouvinte, Em Curso ou Aprovado(a) em alguma disciplina de Didática"
$ws.Cells.Item(25,5).Value = "Projetos pessoais (websites, aplicativos próprios e etc), Contribuições OpenSource"
$ws.Cells.Item(25,6).Value = "Sim"

# Give the three new rows the same visual style as the rest of the table.
$dstA2 = $ws.Range("A23:A25")
$srcA.Copy()
$dstA2.PasteSpecial(-4122)
$dstA2.WrapText = $false

$dstBCF3 = $ws.Range("B23:C25")
$srcB.Copy()
$dstBCF3.PasteSpecial(-4122)
$dstBCF3.WrapText = $false

$dstF3 = $ws.Range("F23:F25")
$srcB.Copy()
$dstF3.PasteSpecial(-4122)
$dstF3.WrapText = $false

$dstDE3 = $ws.Range("D23:E25")
$srcD.Copy()
$dstDE3.PasteSpecial(-4122)
$ws.Range("D23:E25").VerticalAlignment = -4108
